# Add a new bullet paragraph to the "Future Work" slide's content
# placeholder, describing the Twitter-traffic predictor idea.
#
# The new paragraph is typed as two runs (mirroring how PowerPoint splits
# a run when the tail word gets flagged by the spell checker):
#   1) "Create a predictor that leverages Twitter traffic and tweet "
#   2) "classificiation"

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(19)
$sh = $s.Shapes.Item("Content Placeholder 2")
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Start a new paragraph after the existing last bullet ("Analyze links,
# videos, and photos in tweet") and type the first run of text.
$tr.InsertAfter([char]13 + "Create a predictor that leverages Twitter traffic and tweet ") | Out-Null

# Type the (misspelled) second run right after the first one; PowerPoint
# keeps it as its own run.
$tf.TextRange.InsertAfter("classificiation") | Out-Null
